# PAS-6576: fixed VIN upload files for choice and select products,
# changed BI/PD/UM/MP symbols values (swap C <-> A) on rows 2 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (BI_SYMBOL, PD_SYMBOL, UM_SYMBOL, MP_SYMBOL): "A" -> "C"
$ws.Range("AE2:AH2").Value = "C"

# Row 4 (BI_SYMBOL, PD_SYMBOL, UM_SYMBOL, MP_SYMBOL): "C" -> "A"
$ws.Range("AE4:AH4").Value = "A"

# Update the saved view state: scroll the window so column Z is the
# left-most visible column and AI8 becomes the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 26
$win.ScrollRow = 1
$ws.Range("AI8").Select()
